# Deploy the implementation guide.
#
# 1. Flip the ValueSet's published Status from "active" to "draft" and bump
#    the Date metadata row to the new publication timestamp (Metadata sheet).
# 2. Re-assert wrap/alignment formatting on the Metadata + "Include from
#    Ferlab.bio CodeS" sheets so the header/body cell styles carry
#    applyAlignment (alignment was already set but not flagged as applied).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update metadata values -------------------------------------------------
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Re-apply alignment (wrap text; vertical top is already set) on the
# styled ranges, so their xf records gain applyAlignment="true" ------------
# Metadata sheet: header row (A1:B1) + all data rows (A2:B14)
$ws1.Range("A1:B1").WrapText = $true
$ws1.Range("A2:B14").WrapText = $true

# "Include from Ferlab.bio CodeS" sheet: only the cells that actually carry
# data (B1/B2 are genuinely blank/unstyled, so they are left untouched)
$ws2.Range("A1").WrapText = $true
$ws2.Range("A2:A4").WrapText = $true
$ws2.Range("B3:B4").WrapText = $true
